# feat: add 2022-Q1 data
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert a new "2022-Q1" worksheet right before the "总计" sheet,
#    mirroring the layout used by the other quarterly sheets.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$zongji = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($zongji)
$q1.Name = "2022-Q1"

# re-fetch "总计" - Add() can leave stale/shifted sheet references behind
$zongji = $wb.Worksheets.Item("总计")

# Header row + the "A" index column share the exact same layout/style as
# the other quarterly sheets, so copy them (value + format together).
$template.Range("B1:H1").Copy($q1.Range("B1:H1"))
$template.Range("A2").Copy($q1.Range("A2"))

$q1.Cells.Item(2, 2).Value = "'160639"
$q1.Cells.Item(2, 3).Value = "鹏华中证高铁产业指数（LOF）"
$q1.Cells.Item(2, 4).Value = "'0.89"
$q1.Cells.Item(2, 5).Value = "'94.72"
$q1.Cells.Item(2, 6).Value = "'2.62"
$q1.Cells.Item(2, 7).Value = "'0.0233"
$q1.Cells.Item(2, 8).Value = 8

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing quarters down by one row.
# ---------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()

$zongji.Cells.Item(2, 2).Value = "2022-Q1"
$zongji.Cells.Item(2, 3).Value = 1
$zongji.Cells.Item(2, 4).Value = 0.02
$zongji.Range("A3").Copy($zongji.Range("A2"))
$zongji.Cells.Item(2, 1).Value = 0

# Renumber the "A" index column (0-based row counter) for every row that
# shifted down.
$zongji.Cells.Item(3, 1).Value = 1
$zongji.Cells.Item(4, 1).Value = 2
$zongji.Cells.Item(5, 1).Value = 3
$zongji.Cells.Item(6, 1).Value = 4
$zongji.Cells.Item(7, 1).Value = 5
